$d = $word.ActiveDocument

# The target paragraph currently holds ONE run whose text is:
#   "S'il y a un commentaire, il ne doit pas inclure le texte inscrit dans le
#    champ Nom. Sinon, afficher une erreur qui dit : « Le commentaire ne doit
#    pas mentionner le nom de l'utilisateur »"
# The edit splits it into three runs, inserting an extra space between
# "Nom. " and "Sinon":
#   run1: "...champ Nom. "
#   run2: " "
#   run3: "Sinon, afficher...l'utilisateur »"

# Step 1: locate and remove the trailing part of the sentence (from "Sinon"
# onward) - this leaves the run ending in "...Nom. ".
$tail = $d.Content
$tailFound = $tail.Find.Execute("Sinon, afficher une erreur qui dit : « Le commentaire ne doit pas mentionner le nom de l’utilisateur »")
if (-not $tailFound) {
    throw "Could not find the tail sentence to split out"
}
$tail.Delete()

# Step 2: re-find the now-trailing "Nom. " text and re-insert the removed
# tail as two new runs: a standalone space run, then the rest of the
# sentence, appended right after the existing run (still inside the same
# paragraph).
$head = $d.Content
$headFound = $head.Find.Execute("Nom. ")
if (-not $headFound) {
    throw "Could not find 'Nom. ' after trimming the tail"
}

$insertedXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t>Sinon, afficher une erreur qui dit : « Le commentaire ne doit pas mentionner le nom de l’utilisateur »</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$head.InsertXML($insertedXml)
